# Populate the "Expense" sheet's column C (rows 2-11) with the value "sam".
# This reproduces a previously-mis-resolved merge conflict where an extra
# column of placeholder data ("sam") was reintroduced between the
# "Expense Name" (B) and "Expense Category" (D) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expense")

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = "sam"
}

# Match the resulting selection left behind in the sheet (C2:C11, active cell C2).
$ws.Range("C2:C11").Select() | Out-Null
